# LOB1210.xlsx edit
#
# The source table (A:C, rows 1-24) lists course-catalog fields in column A
# with duplicated Portuguese/English values in B/C. The commit removes the
# standalone row 13 (which held the "5840942 - Marco Aurelio Kondracki de
# Alcantara" value with no label in column A) and re-flows the remaining
# content upward - which also drags a handful of adjacent field values one
# slot out of alignment with their labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray, unlabeled row 13. Excel shifts rows 14-24 up to
# become rows 13-23 (and updates the sheet dimension to A1:C23).
$ws.Rows(13).Delete()

# After the shift, a few B/C cells still hold the text that belongs to the
# row above/below; line them back up with their (now relocated) labels.
$ws.Range("B10").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C10").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2020" already exists elsewhere in the sheet as literal text (not a
# date value) - copy it instead of assigning the string directly, so it
# doesn't get auto-converted into a date serial number by the Value setter.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

$ws.Range("B18").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C18").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"

$ws.Range("B19").Value = "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas notas serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Range("C19").Value = "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas notas serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."

$ws.Range("B20").Value = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Range("C20").Value = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."

$ws.Range("B21").Value = "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2."
$ws.Range("C21").Value = "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2."
